$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16 data - copy formatting from row 15 (A15 has the bold/bordered style)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9865283727264408
$ws.Range("D16").Value = 1.428228883633603
$ws.Range("E16").Value = 0.869706806949326
$ws.Range("F16").Value = 0.9865283727264408
$ws.Range("G16").Value = 1.184819198381741
$ws.Range("H16").Value = 0.773072930309377
$ws.Range("I16").Value = 0.898727318825298
$ws.Range("J16").Value = 1.428228883633603
$ws.Range("K16").Value = 1.148967845291465
$ws.Range("L16").Value = 1.067748109008953
$ws.Range("M16").Value = 1.023513918470964

$excel.CutCopyMode = $false
